$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on the D/E columns we are about to rewrite so that
# numeric-looking strings (e.g. "0.600", "68.066.70") keep their original
# text representation instead of being coerced into Excel numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.066.70"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "3.531.70"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "601.28"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "183.38"
$ws.Range("E6").Value = "  +5.51%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "0.141"
$ws.Range("E9").Value = "  +4.75%  "
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "0.442"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "4.142.76"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "32.55"
$ws.Range("E13").Value = "  +11.80%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "68.028.31"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "3.529.07"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "6.42"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "14.94"
$ws.Range("E19").Value = "  +4.64%  "
$ws.Range("D20").Value = "399.59"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "0.548"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "73.62"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "0.0000126"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").Value = "5.70"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "6.34"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "24.15"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "7.52"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("D37").Value = "164.08"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("D39").Value = "0.881"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").Value = "7.22"
$ws.Range("E40").Value = "  +4.47%  "
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "27.14"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").Value = "27.84"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "2.884.87"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").Value = "0.0744"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "354.11"
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").Value = "0.0308"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "1.09"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "34.17"
$ws.Range("E51").Value = "  +1.99%  "

# Reset cell style back to Normal so no stray number-format style lingers
# on cells (keeps styles.xml / cell "s" attributes identical to original).
$textRange.Style = "Normal"
